$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current column B (SFIA Level),
# shifting SFIA Level/Keycode/Description one column to the right.
$ws.Columns.Item(2).Insert()

# New header for the inserted column
$ws.Cells.Item(1, 2).Value = "Skill Description"

# Map each SkillCode (column A) to its full skill description/name
$skillNames = @{
    "Autonomy"   = "Autonomy";
    "Influence"  = "Influence";
    "Complexity" = "Complexity";
    "Knowledge"  = "Knowledge";
    "INOV"       = "Innovation";
    "CNSL"       = "Consultancy";
    "OCDV"       = "Organisational capability development";
    "ORDI"       = "Organisation design and implementation";
    "RLMT"       = "Stakeholder relationship management";
    "KNOW"       = "Knowledge management";
}

$lastRow = 37
for ($r = 2; $r -le $lastRow; $r++) {
    $code = $ws.Cells.Item($r, 1).Value()
    if ($code -ne $null -and $code -ne "") {
        $ws.Cells.Item($r, 2).Value = $skillNames[$code]
    }
}

Write-Output "done"
